# Update "想去人数" (interested-people count) figures in the sheets that
# track event stats: "展览" (sheet1) and "全部类型" (sheet4).
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 545
    $ws.Range("F3").Value = 3508
    $ws.Range("F4").Value = 97
    $ws.Range("F5").Value = 684
}
